{"js": "// This script applies the edit described by the source diff to the\n// document's 9 existing paragraphs and appends one brand-new paragraph:\n//   - paragraphs[0]: date \"03.09.24\" -> \"02.09.24\"\n//   - paragraphs[1]: new title text, plus a manual line break (<w:br/>)\n//     appended at the end of the same run\n//   - paragraphs[2..7]: body text replaced with the new review content\n//   - paragraphs[8]: old arXiv link replaced with a closing sentence\n//   - a new paragraph 9 is appended holding the new arXiv link\nconst newParagraphTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 02.09.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"Transfusion: Predict the Next Token and Diffuse Images with One Multi-Modal Mode\",\n  \"\u05d4\u05d9\u05d5\u05dd \u05e0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05e2\u05dc \u05de\u05d5\u05d3\u05dc \u05de\u05d5\u05dc\u05d8\u05d9\u05de\u05d5\u05d3\u05dc\u05d9 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea. \u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d0\u05d9\u05de\u05e0\u05d5 \u05d1\u05de\u05d0\u05de\u05e8 \u05d9\u05d5\u05d3\u05e2 \u05dc\u05d2\u05e0\u05e8\u05d8 \u05d2\u05dd \u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05d5\u05d2\u05dd \u05d3\u05d0\u05d8\u05d4 \u05d8\u05e7\u05e1\u05d8\u05d5\u05d0\u05dc\u05d9 \u05d5\u05de\u05d4\u05d5\u05d5\u05d4 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d5\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \",\n  \"\u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d6\u05d4 \u05de\u05ea\u05d1\u05d8\u05d0\u05ea \u05d1\u05db\u05da \u05e9\u05d4\u05d9\u05d0 \u05de\u05d2\u05e0\u05e8\u05d8\u05ea \u05d2\u05dd \u05d0\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d8\u05e7\u05e1\u05d8\u05d5\u05d0\u05dc\u05d9 \u05d5\u05d2\u05dd \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05d0\u05e0\u05d5 \u05de\u05d2\u05e0\u05e8\u05d8\u05d9\u05dd \u05d8\u05e7\u05e1\u05d8\u05d9\u05dd, \u05db\u05dc\u05d5\u05de\u05e8 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (\u05e2\u05d1\u05d5\u05e8 \u05ea\u05de\u05d5\u05e0\u05d4 \u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05d8\u05d5\u05e7\u05df \u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9 \u05d0\u05d5 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e4\u05d0\u05e5'). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d2\u05e0\u05e8\u05d8 \u05ea\u05de\u05d5\u05e0\u05d4 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8\u05d4 \u05d4\u05de\u05dc\u05d0 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d2\u05e0\u05e8\u05d8 \u05d0\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8 \u05d8\u05d5\u05e7\u05df \u05d5\u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (next token prediction \u05d0\u05d5 NTP) \u05d5\u05d0\u05d7\u05e8\u05d9 \u05e9\u05d9\u05e1\u05d9\u05d9\u05dd \u05d9\u05d2\u05e0\u05e8\u05d8 \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (\u05d1\u05e6\u05d5\u05e8\u05ea NTP \u05d2\u05dd \u05db\u05df). \u05d6\u05d4 \u05d3\u05d9 \u05e0\u05d7\u05de\u05d3 \u05d4\u05d0\u05de\u05ea.\",\n  \"\u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05db\u05d9\u05dc 7 \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05d6\u05d4 \u05d3\u05d9 \u05e6\u05e0\u05d5\u05e2 \u05dc\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d5\u05d2\u05d5\u05d3\u05dc \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 \u05dc\u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd (\u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d2\u05d3\u05d5\u05dc \u05e9\u05dc stable diffusion \u05de\u05db\u05d9\u05dc \u05d1\u05e2\u05e8\u05da 8B \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd). \u05d0\u05d1\u05dc \u05db\u05d0\u05df \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05e9\u05dc\u05d1 \u05d0\u05ea \u05e9\u05ea\u05d9 \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 (\u05d2\u05e0\u05e8\u05d5\u05d8 \u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05d5\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d8\u05e7\u05e1\u05d8\u05d9\u05dd) \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d3\u05d9 \u05d2\u05d1\u05d5\u05d4\u05d4. \",\n  \"\u05d0\u05d1\u05dc \u05d0\u05d9\u05df \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d6\u05d4? \u05d1\u05d2\u05d3\u05d5\u05dc \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05e7\u05dc\u05d8 \u05e9\u05d4\u05d5\u05d0 \u05e2\u05e8\u05d1\u05d5\u05d1 \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05d5\u05d8\u05e7\u05e1\u05d8 (\u05dc\u05de\u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05e2\u05d5\u05e8\u05d1\u05d1\u05ea \u05e2\u05dd \u05d8\u05e7\u05e1\u05d8). \u05e2\u05dd \u05d4\u05d8\u05e7\u05e1\u05d8 \u05d4\u05db\u05dc \u05e4\u05e9\u05d5\u05d8, \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05d5 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df. \u05dc\u05e4\u05e0\u05d9 \u05db\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d8\u05d5\u05e7\u05df BOI \u05d4\u05de\u05e1\u05de\u05df \u05d0\u05ea \u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d5\u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05d5\u05d6\u05e0\u05d5 \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d8\u05d5\u05e7\u05df EOI \u05dc\u05e1\u05d9\u05de\u05d5\u05df \u05e1\u05d9\u05d5\u05dd \u05d4\u05d6\u05e0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4. \u05db\u05d0\u05de\u05d5\u05e8 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05d6\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd \u05d4\u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d0\u05e6'\u05d9\u05dd \u05dc\u05d0\u05d7\u05e8 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 (\u05e9\u05dc VAE).\",\n  \"\u05d0\u05d9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05d7\u05d9\u05d4 \u05d4\u05d6\u05d5? \u05dc\u05d8\u05e7\u05e1\u05d8 \u05d6\u05d4 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 - \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05d7\u05d6\u05d5\u05ea \u05d8\u05d5\u05e7\u05df \u05d8\u05d5\u05e7\u05df \u05db\u05de\u05d5 \u05d1-LLM \u05e2\u05d1\u05d5\u05e8 \u05de\u05d9\u05dc\u05d5\u05df \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e0\u05ea\u05d5\u05df. \u05e2\u05d1\u05d5\u05e8 \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05dc\u05e4\u05d0\u05e6\u05d9\u05dd, \u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05db\u05dc \u05e4\u05d0\u05e5 \u05d3\u05e8\u05da \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc VAE \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d4 \u05db\u05d8\u05d5\u05e7\u05df. \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05dd \u05de\u05d5\u05e2\u05d1\u05e8\u05d9\u05dd \u05d3\u05e8\u05da \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d0\u05d5 unet \u05dc\u05d4\u05d5\u05e8\u05d3\u05ea \u05de\u05d9\u05de\u05d3. \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05dc\u05d5\u05de\u05d3\u05d9\u05dd \u05dc\u05d4\u05e1\u05d9\u05e8 \u05e8\u05e2\u05e9 \u05de\u05d4\u05d2\u05e8\u05e1\u05d0\u05d5\u05ea \u05d4\u05de\u05d5\u05e8\u05e2\u05e9\u05d5\u05ea \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd. \",\n  \"\u05d1\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d5\u05e6\u05e8 \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05e4\u05d0\u05e5' \u05e4\u05d0\u05e5' \u05de\u05d4\u05e8\u05e2\u05e9 (\u05d0\u05d7\u05e8\u05d9 \u05d4\u05e1\u05e8\u05ea \u05d4\u05e8\u05e2\u05e9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05de\u05d5\u05d6\u05df \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc VAE \u05db\u05d3\u05d9 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05d0\u05e5' \u05e2\u05e6\u05de\u05d5). \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5 \u05dc\u05d9\u05e6\u05d9\u05e8\u05ea \u05ea\u05de\u05d5\u05e0\u05d4 \u05dc\u05d0 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 - \u05e8\u05d5\u05d1 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d9\u05d5\u05e6\u05e8\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05dc\u05d0\u05d4 (\u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc\u05d4). \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d0\u05d7\u05d3 \u05d2\u05d3\u05d5\u05dc!\",\n  \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5 \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\",\n];\nconst newFinalParagraphText = \"https://arxiv.org/pdf/2408.11039\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newParagraphTexts.length) {\n  throw new Error(\n    \"Expected \" + newParagraphTexts.length + \" paragraphs, found \" +\n      paragraphs.items.length\n  );\n}\n\n// Replace the text of every existing paragraph in place so each keeps its\n// original paragraph style (Normal) and run formatting.\nfor (let i = 0; i < newParagraphTexts.length; i++) {\n  paragraphs.items[i].insertText(newParagraphTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// The second paragraph (\"Transfusion: ...\") gets a manual line break\n// appended right after its text -- \"\\v\" (vertical tab) is how Word's\n// manual-line-break character is expressed in a JS string; it serializes\n// as a <w:br/> following the <w:t> inside the same run.\nparagraphs.items[1].insertText(\"\\v\", Word.InsertLocation.end);\nawait context.sync();\n\n// Append a brand-new paragraph with the new arXiv link right after the\n// last paragraph (which now reads \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5 \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\").\nparagraphs.items[newParagraphTexts.length - 1].insertParagraph(\n  newFinalParagraphText,\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# This script applies the edit described by the source diff to the\n# document's 9 existing paragraphs and appends one brand-new paragraph:\n#   - paragraph 1: date \"03.09.24\" -> \"02.09.24\"\n#   - paragraph 2: new title text, plus a manual line break (<w:br/>)\n#     appended at the end of the same run\n#   - paragraphs 3-8: body text replaced with the new review content\n#   - paragraph 9: old arXiv link replaced with a closing sentence\n#   - a new paragraph 10 is appended holding the new arXiv link\n\n$d = $word.ActiveDocument\n\n$newParagraphTexts = @(\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 02.09.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"Transfusion: Predict the Next Token and Diffuse Images with One Multi-Modal Mode\",\n  \"\u05d4\u05d9\u05d5\u05dd \u05e0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05e2\u05dc \u05de\u05d5\u05d3\u05dc \u05de\u05d5\u05dc\u05d8\u05d9\u05de\u05d5\u05d3\u05dc\u05d9 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea. \u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d0\u05d9\u05de\u05e0\u05d5 \u05d1\u05de\u05d0\u05de\u05e8 \u05d9\u05d5\u05d3\u05e2 \u05dc\u05d2\u05e0\u05e8\u05d8 \u05d2\u05dd \u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05d5\u05d2\u05dd \u05d3\u05d0\u05d8\u05d4 \u05d8\u05e7\u05e1\u05d8\u05d5\u05d0\u05dc\u05d9 \u05d5\u05de\u05d4\u05d5\u05d5\u05d4 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d5\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \",\n  \"\u05d4\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d6\u05d4 \u05de\u05ea\u05d1\u05d8\u05d0\u05ea \u05d1\u05db\u05da \u05e9\u05d4\u05d9\u05d0 \u05de\u05d2\u05e0\u05e8\u05d8\u05ea \u05d2\u05dd \u05d0\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d8\u05e7\u05e1\u05d8\u05d5\u05d0\u05dc\u05d9 \u05d5\u05d2\u05dd \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05d0\u05e0\u05d5 \u05de\u05d2\u05e0\u05e8\u05d8\u05d9\u05dd \u05d8\u05e7\u05e1\u05d8\u05d9\u05dd, \u05db\u05dc\u05d5\u05de\u05e8 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (\u05e2\u05d1\u05d5\u05e8 \u05ea\u05de\u05d5\u05e0\u05d4 \u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05d8\u05d5\u05e7\u05df \u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9 \u05d0\u05d5 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05e4\u05d0\u05e5'). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d2\u05e0\u05e8\u05d8 \u05ea\u05de\u05d5\u05e0\u05d4 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05ea\u05d9\u05d0\u05d5\u05e8\u05d4 \u05d4\u05de\u05dc\u05d0 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d2\u05e0\u05e8\u05d8 \u05d0\u05ea \u05d4\u05ea\u05d9\u05d0\u05d5\u05e8 \u05d8\u05d5\u05e7\u05df \u05d5\u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (next token prediction \u05d0\u05d5 NTP) \u05d5\u05d0\u05d7\u05e8\u05d9 \u05e9\u05d9\u05e1\u05d9\u05d9\u05dd \u05d9\u05d2\u05e0\u05e8\u05d8 \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df (\u05d1\u05e6\u05d5\u05e8\u05ea NTP \u05d2\u05dd \u05db\u05df). \u05d6\u05d4 \u05d3\u05d9 \u05e0\u05d7\u05de\u05d3 \u05d4\u05d0\u05de\u05ea.\",\n  \"\u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05db\u05d9\u05dc 7 \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05d6\u05d4 \u05d3\u05d9 \u05e6\u05e0\u05d5\u05e2 \u05dc\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d5\u05d2\u05d5\u05d3\u05dc \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 \u05dc\u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd (\u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d2\u05d3\u05d5\u05dc \u05e9\u05dc stable diffusion \u05de\u05db\u05d9\u05dc \u05d1\u05e2\u05e8\u05da 8B \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd). \u05d0\u05d1\u05dc \u05db\u05d0\u05df \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05e9\u05dc\u05d1 \u05d0\u05ea \u05e9\u05ea\u05d9 \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 (\u05d2\u05e0\u05e8\u05d5\u05d8 \u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05d5\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d8\u05e7\u05e1\u05d8\u05d9\u05dd) \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d3\u05d9 \u05d2\u05d1\u05d5\u05d4\u05d4. \",\n  \"\u05d0\u05d1\u05dc \u05d0\u05d9\u05df \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d6\u05d4? \u05d1\u05d2\u05d3\u05d5\u05dc \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05e7\u05dc\u05d8 \u05e9\u05d4\u05d5\u05d0 \u05e2\u05e8\u05d1\u05d5\u05d1 \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05d5\u05d8\u05e7\u05e1\u05d8 (\u05dc\u05de\u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05e2\u05d5\u05e8\u05d1\u05d1\u05ea \u05e2\u05dd \u05d8\u05e7\u05e1\u05d8). \u05e2\u05dd \u05d4\u05d8\u05e7\u05e1\u05d8 \u05d4\u05db\u05dc \u05e4\u05e9\u05d5\u05d8, \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05d5 \u05d8\u05d5\u05e7\u05df \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d5\u05e7\u05df. \u05dc\u05e4\u05e0\u05d9 \u05db\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d8\u05d5\u05e7\u05df BOI \u05d4\u05de\u05e1\u05de\u05df \u05d0\u05ea \u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d5\u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05d5\u05d6\u05e0\u05d5 \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d8\u05d5\u05e7\u05df EOI \u05dc\u05e1\u05d9\u05de\u05d5\u05df \u05e1\u05d9\u05d5\u05dd \u05d4\u05d6\u05e0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4. \u05db\u05d0\u05de\u05d5\u05e8 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05d6\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd \u05d4\u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d0\u05e6'\u05d9\u05dd \u05dc\u05d0\u05d7\u05e8 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 (\u05e9\u05dc VAE).\",\n  \"\u05d0\u05d9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05d7\u05d9\u05d4 \u05d4\u05d6\u05d5? \u05dc\u05d8\u05e7\u05e1\u05d8 \u05d6\u05d4 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 - \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05d7\u05d6\u05d5\u05ea \u05d8\u05d5\u05e7\u05df \u05d8\u05d5\u05e7\u05df \u05db\u05de\u05d5 \u05d1-LLM \u05e2\u05d1\u05d5\u05e8 \u05de\u05d9\u05dc\u05d5\u05df \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e0\u05ea\u05d5\u05df. \u05e2\u05d1\u05d5\u05e8 \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05dc\u05e4\u05d0\u05e6\u05d9\u05dd, \u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05db\u05dc \u05e4\u05d0\u05e5 \u05d3\u05e8\u05da \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc VAE \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d4 \u05db\u05d8\u05d5\u05e7\u05df. \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05dd \u05de\u05d5\u05e2\u05d1\u05e8\u05d9\u05dd \u05d3\u05e8\u05da \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d0\u05d5 unet \u05dc\u05d4\u05d5\u05e8\u05d3\u05ea \u05de\u05d9\u05de\u05d3. \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05dc\u05d5\u05de\u05d3\u05d9\u05dd \u05dc\u05d4\u05e1\u05d9\u05e8 \u05e8\u05e2\u05e9 \u05de\u05d4\u05d2\u05e8\u05e1\u05d0\u05d5\u05ea \u05d4\u05de\u05d5\u05e8\u05e2\u05e9\u05d5\u05ea \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d5\u05d9\u05d6\u05d5\u05d0\u05dc\u05d9\u05d9\u05dd. \",\n  \"\u05d1\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d5\u05e6\u05e8 \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05e4\u05d0\u05e5' \u05e4\u05d0\u05e5' \u05de\u05d4\u05e8\u05e2\u05e9 (\u05d0\u05d7\u05e8\u05d9 \u05d4\u05e1\u05e8\u05ea \u05d4\u05e8\u05e2\u05e9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05de\u05d5\u05d6\u05df \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc VAE \u05db\u05d3\u05d9 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05d0\u05e5' \u05e2\u05e6\u05de\u05d5). \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5 \u05dc\u05d9\u05e6\u05d9\u05e8\u05ea \u05ea\u05de\u05d5\u05e0\u05d4 \u05dc\u05d0 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 - \u05e8\u05d5\u05d1 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05d9\u05d5\u05e6\u05e8\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05dc\u05d0\u05d4 (\u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc\u05d4). \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d0\u05d7\u05d3 \u05d2\u05d3\u05d5\u05dc!\",\n  \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5 \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\"\n)\n$newFinalParagraphText = \"https://arxiv.org/pdf/2408.11039\"\n\nif ($d.Paragraphs.Count -ne $newParagraphTexts.Count) {\n  throw \"Expected $($newParagraphTexts.Count) paragraphs, found $($d.Paragraphs.Count)\"\n}\n\n# Replace the text of every existing paragraph in place so each keeps its\n# original paragraph style (Normal) and run formatting.\nfor ($i = 0; $i -lt $newParagraphTexts.Count; $i++) {\n  $d.Paragraphs($i + 1).Range.Text = $newParagraphTexts[$i]\n}\n\n# The second paragraph (\"Transfusion: ...\") gets a manual line break\n# appended right after its text -- [char]11 is Word's manual-line-break\n# character (vertical tab); it serializes as a <w:br/> following the\n# <w:t> inside the same run.\n$d.Paragraphs(2).Range.InsertAfter([char]11)\n\n# Append a brand-new paragraph with the new arXiv link right after the\n# last paragraph (which now reads \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5 \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\").\n$lastIndex = $newParagraphTexts.Count\n$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()\n$d.Paragraphs($lastIndex + 1).Range.Text = $newFinalParagraphText\n"}
